# Add a "Save" column (column H) to the s_vals sheet, mirroring the
# header formatting used by the existing "sum" column (G) and filling
# in the value for the single data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the neighboring header cell (G1, bold + border +
# centered/top aligned) onto the new header cell H1, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Data cell for the new column.
$ws.Range("H2").Value = 1
